# bug fix in Eduati data files
# HCT116_noCTRL_meas.xlsx: Sheet1 had 43 extra stray rows (45:87) that only
# carried a leftover index in column A with no real data alongside them.
# Trim Sheet1 back down to the real data range (A1:N44, same shape as
# Sheet2/Sheet3), and restore the view state to what was captured when the
# file was actually re-saved: Sheet1 becomes the active/selected tab (so the
# workbook no longer opens on Sheet3), and the stray "tabSelected" marker on
# Sheet3 is cleared.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Sheet1: drop the stray tail rows (45 through 87) ---------------------
$ws1.Activate()
$lastRow = $ws1.UsedRange.Rows.Count
if ($lastRow -ge 45) {
    $staleRows = $ws1.Range($ws1.Cells.Item(45, 1), $ws1.Cells.Item($lastRow, 1))
    $staleRows.EntireRow.Delete()
}

# --- Restore per-sheet selections to match the real saved file ------------
$ws2.Range("A2:N44").Select()
$ws3.Range("A2:N44").Select()

# Sheet1 ends up active/selected (tabSelected) and on cell H40, matching the
# post-edit workbook (workbook-level activeTab moves off Sheet3 back to
# Sheet1 as a result). Activating it last also clears "tabSelected" back off
# of Sheet2/Sheet3.
$ws1.Activate()
$ws1.Range("H40").Select()
